# Apply weekly update: insert two new data rows (new rows 163 and 164)
# before the former row 163, shifting the rest of the data block down by
# two rows (old rows 163-269 become new rows 165-271).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 163:164 - this shifts existing rows 163-269
# down to 165-271 and expands the used range/dimension automatically.
$ws.Rows("163:164").Insert()

# --- New row 163 ---
$ws.Range("A163").Value = 8
$ws.Range("B163").Value = "Terminal La Palmera de La Serena"
$ws.Range("C163").Value = "Coquimbo"
$ws.Range("D163").Value = 45126
$ws.Range("E163").Value = 4
$ws.Range("F163").Value = 100112001
$ws.Range("G163").Value = "Berenjena"
$ws.Range("H163").Value = "Sin especificar"
$ws.Range("I163").Value = "Primera"
$ws.Range("J163").Value = 440
$ws.Range("K163").Value = 8000
$ws.Range("L163").Value = 9000
$ws.Range("M163").Value = 8500
$ws.Range("N163").Value = "$/caja 50 unidades"
$ws.Range("O163").Value = "Región de Arica y Parinacota"
$ws.Range("P163").Value = 170
$ws.Range("Q163").Value = 50
$ws.Range("R163").Value = "Hortaliza"

# --- New row 164 ---
$ws.Range("A164").Value = 8
$ws.Range("B164").Value = "Terminal La Palmera de La Serena"
$ws.Range("C164").Value = "Coquimbo"
$ws.Range("D164").Value = 45126
$ws.Range("E164").Value = 4
$ws.Range("F164").Value = 100112001
$ws.Range("G164").Value = "Berenjena"
$ws.Range("H164").Value = "Sin especificar"
$ws.Range("I164").Value = "Primera"
$ws.Range("J164").Value = 360
$ws.Range("K164").Value = 8000
$ws.Range("L164").Value = 9000
$ws.Range("M164").Value = 8500
$ws.Range("N164").Value = "$/caja 50 unidades"
$ws.Range("O164").Value = "Región de Arica y Parinacota"
$ws.Range("P164").Value = 170
$ws.Range("Q164").Value = 50
$ws.Range("R164").Value = "Hortaliza"

# Make sure the date cells keep the date number format used throughout
# column D (same as the cell above / below them).
$ws.Range("D163:D164").NumberFormat = $ws.Range("D165").NumberFormat
